$d = $word.ActiveDocument

$replacements = @(
    @("99×44=", "17×65="),
    @("78×35=", "12×99="),
    @("34×78=", "66×81="),
    @("99×69=", "71×11="),
    @("67×55=", "43×15="),
    @("42×25=", "57×82="),
    @("13×59=", "35×58="),
    @("57×58=", "39×72="),
    @("33×20=", "38×65="),
    @("87×83=", "51×86="),
    @("72×92=", "44×20="),
    @("24×75=", "33×29="),
    @("18×56=", "47×17="),
    @("46×49=", "28×54="),
    @("66×11=", "86×58="),
    @("15×48=", "25×60="),
    @("60×48=", "11×97="),
    @("98×45=", "84×77="),
    @("47×86=", "82×80="),
    @("57×66=", "29×70="),
    @("47×76=", "21×78="),
    @("64×11=", "49×72="),
    @("83×36=", "70×53="),
    @("48×47=", "50×68="),
    @("67×16=", "55×86=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
